# Update "想去人数" (F column) values on several sheets to reflect the
# latest scrape output (gh-pages rebuild at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "展览"
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 668
$ws.Range("F4").Value  = 7174
$ws.Range("F6").Value  = 13324
$ws.Range("F7").Value  = 13477
$ws.Range("F9").Value  = 1335
$ws.Range("F10").Value = 5650
$ws.Range("F19").Value = 1105
$ws.Range("F23").Value = 2208
$ws.Range("F26").Value = 3130
$ws.Range("F27").Value = 280
$ws.Range("F29").Value = 33
$ws.Range("F36").Value = 4438
$ws.Range("F37").Value = 4568
$ws.Range("F45").Value = 314
$ws.Range("F49").Value = 230

# -----------------------------------------------------------------
# Sheet "本地生活"
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 7154
$ws.Range("F3").Value = 161
$ws.Range("F4").Value = 479

# -----------------------------------------------------------------
# Sheet "全部类型"
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 668
$ws.Range("F5").Value  = 161
$ws.Range("F6").Value  = 479
$ws.Range("F8").Value  = 13324
$ws.Range("F9").Value  = 13477
$ws.Range("F12").Value = 5650
$ws.Range("F20").Value = 1105
$ws.Range("F25").Value = 3130
$ws.Range("F27").Value = 280
$ws.Range("F29").Value = 33
$ws.Range("F37").Value = 4438
$ws.Range("F38").Value = 4568
$ws.Range("F46").Value = 314

$wb.Save()
